$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = -0.8319244784142803
$ws.Range("J2").Value = 0.2980185349714233
$ws.Range("K2").Value = 0.02880460297801311
$ws.Range("L2").Value = 2.409774392370587

$ws.Range("I19").Value = -1.064347500819291
$ws.Range("J19").Value = 0.3685536603139039
$ws.Range("K19").Value = 0.3361587556007758
$ws.Range("L19").Value = 2.194017056983228
